$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update correct error codes in test-challenges
$ws.Range("D22").Value = 1
$ws.Range("D33").Value = 1
$ws.Range("D36").Value = 2

# Update active selection to D22
$ws.Range("D22").Select()
